# Update cryptocurrency price and volume figures (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.226.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.907.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5251"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3784"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07275"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8993"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07694"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.910.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008665"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.288.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.089"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.155.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.455"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.327"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.738"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.977"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.821"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09243"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8180"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05080"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  +7.87%  "
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.314"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.606"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5676"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01993"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.077"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.018"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.653"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.92%  "
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4849"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.621"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.42%  "
